# Add two new shortcut rows ("Store Zoom" / ZS and "Restore Zoom" / ZR)
# right above the existing "Flip X-Y Axis" row (old row 104), shifting
# everything below down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the old row 104 ("Flip X-Y Axis"); this
# pushes the old rows 104-121 down to 106-123 and copies formatting
# (styles/row height) from the row below, matching the surrounding rows.
$ws.Rows("104:105").Insert()
$ws.Rows("104:105").RowHeight = 17

# Fill in the new shortcut entries.
$ws.Range("A104").Value = "Store Zoom"
$ws.Range("B104").Value = "ZS"
$ws.Range("A105").Value = "Restore Zoom"
$ws.Range("B105").Value = "ZR"

# Move the active selection to reflect where the user ended up editing.
$ws.Range("B106").Select()
